$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data refresh (16th May Refresh): append new machine rows for regcntr_id 10005
$newRows = @(
    @(10005, 110033, 10005, "eng", $true, "superadmin", "now()"),
    @(10005, 110034, 10005, "eng", $true, "superadmin", "now()"),
    @(10005, 110035, 10005, "eng", $true, "superadmin", "now()")
)

$startRow = 34
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}

# Mimic the post-entry selection left behind after entering data (selecting
# the remainder of the rows below the table, as seen in the saved file)
$ws.Range("A37:XFD1048576").Select() | Out-Null
